# Added a '10' in sheet 3.
$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("C10").Value = 10

$ws3.Activate()
$ws3.Range("C10").Select()
